$d = $word.ActiveDocument

# --- Part 1: the "又因为...作为自己的一个私有属性。" paragraph ---------------
# Originally every run (and the paragraph mark) carries <w:strike/>.
# The edit removes the strikethrough and turns the text red instead; in
# addition the three runs "Execute" / "和" / "Commit" are collapsed into a
# single run reading "时间更新".

$rng = $d.Content

# "又因为" - format only, no text change
$rng.Find.Execute("又因为", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0) | Out-Null
$rng.Font.StrikeThrough = $false
$rng.Font.Color = 255

# "Execute和Commit" -> "时间更新" (merges the 3 runs into 1) + same format
$rng.Find.Execute("Execute和Commit", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "时间更新", 1) | Out-Null
$rng.Font.StrikeThrough = $false
$rng.Font.Color = 255

# remaining runs - format only, no text change, in document order
$remainingRuns = @(
    "是对",
    "subprogram",
    "执行的操作，而我们需要将预估时间更新到",
    "sub template",
    "上，所以",
    "subprogram",
    "需要将",
    "sub template",
    "作为自己的一个私有属性。"
)

foreach ($t in $remainingRuns) {
    $rng.Find.Execute($t, $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0) | Out-Null
    $rng.Font.StrikeThrough = $false
    $rng.Font.Color = 255
}

# Paragraph mark (pPr/rPr) itself also switches strike -> color.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "又因为*") {
        $p.Range.ParagraphFormat.Application | Out-Null
        $pMarkRange = $p.Range
        $pMarkRange.Start = $pMarkRange.End - 1
        $pMarkRange.Font.StrikeThrough = $false
        $pMarkRange.Font.Color = 255
        break
    }
}
